{"js": "// Replace the date line and each \"AxB=\" multiplication prompt in the\n// document with its updated value. Every source string in this document\n// is unique, so an exact, case-sensitive search+replace per pair is\n// unambiguous and preserves the existing run formatting (font, size).\nconst replacements = [\n  [\"2024-06-26 Wednesday\", \"2024-06-27 Thursday\"],\n  [\"75\u00d719=\", \"25\u00d723=\"],\n  [\"63\u00d733=\", \"65\u00d787=\"],\n  [\"42\u00d774=\", \"54\u00d742=\"],\n  [\"48\u00d756=\", \"58\u00d794=\"],\n  [\"22\u00d714=\", \"57\u00d773=\"],\n  [\"23\u00d781=\", \"41\u00d747=\"],\n  [\"47\u00d745=\", \"29\u00d752=\"],\n  [\"61\u00d716=\", \"95\u00d721=\"],\n  [\"95\u00d745=\", \"86\u00d792=\"],\n  [\"29\u00d792=\", \"84\u00d774=\"],\n  [\"23\u00d714=\", \"65\u00d773=\"],\n  [\"73\u00d790=\", \"27\u00d722=\"],\n  [\"91\u00d785=\", \"33\u00d789=\"],\n  [\"72\u00d714=\", \"13\u00d750=\"],\n  [\"86\u00d783=\", \"87\u00d735=\"],\n  [\"80\u00d772=\", \"38\u00d778=\"],\n  [\"69\u00d791=\", \"18\u00d716=\"],\n  [\"98\u00d798=\", \"71\u00d782=\"],\n  [\"61\u00d738=\", \"89\u00d740=\"],\n  [\"36\u00d760=\", \"35\u00d742=\"],\n  [\"43\u00d799=\", \"65\u00d780=\"],\n  [\"73\u00d784=\", \"12\u00d735=\"],\n  [\"65\u00d737=\", \"33\u00d751=\"],\n  [\"14\u00d761=\", \"76\u00d771=\"],\n  [\"56\u00d712=\", \"62\u00d726=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"AxB=\" multiplication prompt in the\n# document with its updated value. Every source string in this document\n# is unique, so a plain Find/Replace per pair is unambiguous and leaves\n# the existing run formatting (font, size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-06-26 Wednesday\", \"2024-06-27 Thursday\"),\n    @(\"75\u00d719=\", \"25\u00d723=\"),\n    @(\"63\u00d733=\", \"65\u00d787=\"),\n    @(\"42\u00d774=\", \"54\u00d742=\"),\n    @(\"48\u00d756=\", \"58\u00d794=\"),\n    @(\"22\u00d714=\", \"57\u00d773=\"),\n    @(\"23\u00d781=\", \"41\u00d747=\"),\n    @(\"47\u00d745=\", \"29\u00d752=\"),\n    @(\"61\u00d716=\", \"95\u00d721=\"),\n    @(\"95\u00d745=\", \"86\u00d792=\"),\n    @(\"29\u00d792=\", \"84\u00d774=\"),\n    @(\"23\u00d714=\", \"65\u00d773=\"),\n    @(\"73\u00d790=\", \"27\u00d722=\"),\n    @(\"91\u00d785=\", \"33\u00d789=\"),\n    @(\"72\u00d714=\", \"13\u00d750=\"),\n    @(\"86\u00d783=\", \"87\u00d735=\"),\n    @(\"80\u00d772=\", \"38\u00d778=\"),\n    @(\"69\u00d791=\", \"18\u00d716=\"),\n    @(\"98\u00d798=\", \"71\u00d782=\"),\n    @(\"61\u00d738=\", \"89\u00d740=\"),\n    @(\"36\u00d760=\", \"35\u00d742=\"),\n    @(\"43\u00d799=\", \"65\u00d780=\"),\n    @(\"73\u00d784=\", \"12\u00d735=\"),\n    @(\"65\u00d737=\", \"33\u00d751=\"),\n    @(\"14\u00d761=\", \"76\u00d771=\"),\n    @(\"56\u00d712=\", \"62\u00d726=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
